# Auto-generated Excel COM-interop script to apply the diff to Mateus_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6027.731
$ws.Range("I74").Value = 3540.3333
$ws.Range("K74").Value = 3540.3333
$ws.Range("M74").Value = -2604.3333

$ws.Range("H77").Value = 6027.731
$ws.Range("I77").Value = 3540.3333
$ws.Range("K77").Value = 17701.6665
$ws.Range("M77").Value = -13021.6665

$ws.Range("H92").Value = 613.7143
$ws.Range("I92").Value = 707.6667
$ws.Range("J92").Value = 50
$ws.Range("K92").Value = 707.6667
$ws.Range("L92").Value = 50
$ws.Range("M92").Value = 540.3333
$ws.Range("N92").Value = -2546

$ws.Range("H101").Value = 1178
$ws.Range("I101").Value = 913.6
$ws.Range("J101").Value = 2500
$ws.Range("K101").Value = 2740.8
$ws.Range("L101").Value = 7500
$ws.Range("M101").Value = -1118.8
$ws.Range("N101").Value = -10744

$ws.Range("H106").Value = 2967.9
$ws.Range("I106").Value = 2459.875
$ws.Range("K106").Value = 2459.875
$ws.Range("M106").Value = -1828.875

$ws.Range("H132").Value = 5466.16
$ws.Range("I132").Value = 6291.55
$ws.Range("K132").Value = 18874.65
$ws.Range("M132").Value = -16344.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4848.5
$ws.Range("I2").Value = 5154.9165
$ws.Range("J2").Value = 4480.8
$ws.Range("K2").Value = 5154.9165
$ws.Range("L2").Value = 4480.8
$ws.Range("M2").Value = -5041.9165
$ws.Range("N2").Value = -4706.8

$ws.Range("H97").Value = 742.6667
$ws.Range("I97").Value = 597.7143
$ws.Range("J97").Value = 1250
$ws.Range("K97").Value = 597.7143
$ws.Range("L97").Value = 1250
$ws.Range("M97").Value = -101.7143
$ws.Range("N97").Value = -2242

$ws.Range("H102").Value = 4779.6
$ws.Range("I102").Value = 3724.5
$ws.Range("K102").Value = 3724.5
$ws.Range("M102").Value = -2102.5

$ws.Range("H110").Value = 3933.6155
$ws.Range("I110").Value = 2968.95
$ws.Range("K110").Value = 2968.95
$ws.Range("M110").Value = -923.9499999999998

$ws.Range("H116").Value = 4848.5
$ws.Range("I116").Value = 5154.9165
$ws.Range("J116").Value = 4480.8
$ws.Range("K116").Value = 5154.9165
$ws.Range("L116").Value = 4480.8
$ws.Range("M116").Value = -2860.9165
$ws.Range("N116").Value = -9068.799999999999

$ws.Range("H132").Value = 5615.6113
$ws.Range("I132").Value = 4621.9614
$ws.Range("K132").Value = 13865.8842
$ws.Range("M132").Value = -11335.8842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4848.5
$ws.Range("I3").Value = 5154.9165
$ws.Range("J3").Value = 4480.8
$ws.Range("K3").Value = 5154.9165
$ws.Range("L3").Value = 4480.8
$ws.Range("M3").Value = -5040.9165
$ws.Range("N3").Value = -4708.8

$ws.Range("H82").Value = 38900
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0

$ws.Range("H85").Value = 38900
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0

$ws.Range("H86").Value = 9092350
$ws.Range("I86").Value = 1515.625
$ws.Range("J86").Value = 33334576
$ws.Range("K86").Value = 1515.625
$ws.Range("L86").Value = 33334576
$ws.Range("M86").Value = -392.625
$ws.Range("N86").Value = -33336822

$ws.Range("H89").Value = 9092350
$ws.Range("I89").Value = 1515.625
$ws.Range("J89").Value = 33334576
$ws.Range("K89").Value = 7578.125
$ws.Range("L89").Value = 166672880
$ws.Range("M89").Value = -1962.125
$ws.Range("N89").Value = -166684112

$ws.Range("H94").Value = 1983.6666
$ws.Range("I94").Value = 1652.0625
$ws.Range("J94").Value = 3044.8
$ws.Range("K94").Value = 1652.0625
$ws.Range("L94").Value = 3044.8
$ws.Range("M94").Value = -1201.0625
$ws.Range("N94").Value = -3946.8

$ws.Range("H99").Value = 3691.9512
$ws.Range("I99").Value = 2485.926
$ws.Range("J99").Value = 6017.857
$ws.Range("K99").Value = 2485.926
$ws.Range("L99").Value = 6017.857
$ws.Range("M99").Value = -987.9259999999999
$ws.Range("N99").Value = -9013.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3408.5789
$ws.Range("J16").Value = 3858.818
$ws.Range("L16").Value = 3858.818
$ws.Range("N16").Value = -4432.818

$ws.Range("H22").Value = 3742.25
$ws.Range("I22").Value = 4221.4287
$ws.Range("K22").Value = 4221.4287
$ws.Range("M22").Value = -3871.4287

$ws.Range("H31").Value = 4959.029
$ws.Range("J31").Value = 9831.444
$ws.Range("L31").Value = 9831.444
$ws.Range("N31").Value = -10421.444

$ws.Range("H34").Value = 4959.029
$ws.Range("J34").Value = 9831.444
$ws.Range("L34").Value = 9831.444
$ws.Range("N34").Value = -10235.444

$ws.Range("H51").Value = 30029.666
$ws.Range("I51").Value = 30029
$ws.Range("K51").Value = 30029
$ws.Range("M51").Value = -29293

$ws.Range("H61").Value = 30029.666
$ws.Range("I61").Value = 30029
$ws.Range("K61").Value = 30029
$ws.Range("M61").Value = -29681

$ws.Range("H113").Value = 3408.5789
$ws.Range("J113").Value = 3858.818
$ws.Range("L113").Value = 3858.818
$ws.Range("N113").Value = -8198.817999999999

$ws.Range("H122").Value = 2522.8064
$ws.Range("I122").Value = 2356.9
$ws.Range("K122").Value = 7070.700000000001
$ws.Range("M122").Value = -4620.700000000001

$ws.Range("H132").Value = 1839.258
$ws.Range("I132").Value = 1221.36
$ws.Range("K132").Value = 3664.08
$ws.Range("M132").Value = -1134.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 331011.8
$ws.Range("I128").Value = 331011.8
$ws.Range("K128").Value = 993035.3999999999
$ws.Range("M128").Value = -988055.3999999999

$ws.Range("H140").Value = 523088.7
$ws.Range("I140").Value = 1223.2222
$ws.Range("J140").Value = 950069.5600000001
$ws.Range("K140").Value = 3669.6666
$ws.Range("L140").Value = 2850208.68
$ws.Range("M140").Value = 1510.3334
$ws.Range("N140").Value = -2860568.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3569.25
$ws.Range("I80").Value = 2166.6667
$ws.Range("K80").Value = 2166.6667
$ws.Range("M80").Value = -1168.6667

$ws.Range("H83").Value = 3569.25
$ws.Range("I83").Value = 2166.6667
$ws.Range("K83").Value = 10833.3335
$ws.Range("M83").Value = -5841.333500000001

$ws.Range("H132").Value = 5257.5
$ws.Range("I132").Value = 4931.1924
$ws.Range("K132").Value = 14793.5772
$ws.Range("M132").Value = -12263.5772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2716.8125
$ws.Range("I40").Value = 2497.923
$ws.Range("K40").Value = 2497.923
$ws.Range("M40").Value = -2361.923

$ws.Range("H82").Value = 4956.35
$ws.Range("I82").Value = 3715.4443
$ws.Range("K82").Value = 3715.4443
$ws.Range("M82").Value = -3354.4443

$ws.Range("H85").Value = 4956.35
$ws.Range("I85").Value = 3715.4443
$ws.Range("K85").Value = 3715.4443
$ws.Range("M85").Value = -2467.4443

$ws.Range("H93").Value = 13281.588
$ws.Range("I93").Value = 856.2857
$ws.Range("J93").Value = 71266.336
$ws.Range("K93").Value = 856.2857
$ws.Range("L93").Value = 71266.336
$ws.Range("M93").Value = 391.7143
$ws.Range("N93").Value = -73762.336

$ws.Range("H100").Value = 1616223.8
$ws.Range("I100").Value = 2002917.5
$ws.Range("K100").Value = 2002917.5
$ws.Range("M100").Value = -2002376.5

$ws.Range("H140").Value = 90187.336
$ws.Range("J140").Value = 90187.336
$ws.Range("L140").Value = 90187.336
$ws.Range("N140").Value = -100547.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2817.6924
$ws.Range("I122").Value = 2158.75
$ws.Range("K122").Value = 6476.25
$ws.Range("M122").Value = -4026.25

$ws.Range("H132").Value = 4846.073
$ws.Range("I132").Value = 4274
$ws.Range("K132").Value = 12822
$ws.Range("M132").Value = -10292

# Rows where LeveProfitNQ (M) became 0 and the cell is cleared entirely (matches removal of <c> element in diff)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M82").ClearContents()
$ws.Range("M85").ClearContents()
